$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (this shifts existing rows 85:101 down to 86:102)
$ws.Rows.Item(85).Insert()

# Populate the new row 85 with a new weekly price observation.
# Non price/date columns mirror the record that used to sit at row 85
# (same market/region/category/quality/unit), matching the target diff.
$ws.Cells.Item(85, 1).Value = 11
$ws.Cells.Item(85, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(85, 3).Value = "Bíobío"
$ws.Cells.Item(85, 4).Value = 45009
$ws.Cells.Item(85, 5).Value = 8
$ws.Cells.Item(85, 6).Value = 100112037
$ws.Cells.Item(85, 7).Value = "Cebollín"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 100
$ws.Cells.Item(85, 11).Value = 4500
$ws.Cells.Item(85, 12).Value = 5000
$ws.Cells.Item(85, 13).Value = 4750
$ws.Cells.Item(85, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(85, 15).Value = "Región Metropolitana"
$ws.Cells.Item(85, 16).Value = 132
$ws.Cells.Item(85, 17).Value = 36
$ws.Cells.Item(85, 18).Value = "Hortaliza"
